# Fruta / hortaliza, semanal
# Insert two new weekly-report rows ("Early Majestic") right after the
# existing row 490 ("Elegant Lady" / Segunda), pushing the rows that used
# to be 491-509 down to 493-511.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 491 (Excel shifts 491..509 down to 493..511,
# inheriting the formatting -- including the date-number-format style on
# column D -- from the row immediately above, same as a manual Excel
# "Insert Copied/Sheet Rows" would do).
$ws.Rows.Item(491).Resize(2).Insert()

# Row 491 - Femacal de La Calera / Coquimbo / Early Majestic / Primera
$ws.Range("A491").Value = 3
$ws.Range("B491").Value = 'Femacal de La Calera'
$ws.Range("C491").Value = 'Coquimbo'
$ws.Range("D491").Value = 44509
$ws.Range("E491").Value = 5
$ws.Range("F491").Value = 'Fruta'
$ws.Range("G491").Value = 100103
$ws.Range("H491").Value = 'Frutos de hueso (carozo)'
$ws.Range("I491").Value = 100103004
$ws.Range("J491").Value = 'Durazno'
$ws.Range("K491").Value = 'Early Majestic'
$ws.Range("L491").Value = 'Primera'
$ws.Range("M491").Value = 40
$ws.Range("N491").Value = 13000
$ws.Range("O491").Value = 13000
$ws.Range("P491").Value = 13000
$ws.Range("Q491").Value = '$/bandeja 10 kilos granel'
$ws.Range("R491").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S491").Value = 1300
$ws.Range("T491").Value = 10

# Row 492 - Femacal de La Calera / Coquimbo / Early Majestic / Segunda
$ws.Range("A492").Value = 3
$ws.Range("B492").Value = 'Femacal de La Calera'
$ws.Range("C492").Value = 'Coquimbo'
$ws.Range("D492").Value = 44509
$ws.Range("E492").Value = 5
$ws.Range("F492").Value = 'Fruta'
$ws.Range("G492").Value = 100103
$ws.Range("H492").Value = 'Frutos de hueso (carozo)'
$ws.Range("I492").Value = 100103004
$ws.Range("J492").Value = 'Durazno'
$ws.Range("K492").Value = 'Early Majestic'
$ws.Range("L492").Value = 'Segunda'
$ws.Range("M492").Value = 36
$ws.Range("N492").Value = 10000
$ws.Range("O492").Value = 10000
$ws.Range("P492").Value = 10000
$ws.Range("Q492").Value = '$/bandeja 10 kilos granel'
$ws.Range("R492").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S492").Value = 1000
$ws.Range("T492").Value = 10

# Make sure column D on the two new rows keeps the date-serial number
# format used by the rest of the column (style index 2 in styles.xml).
$ws.Range("D491:D492").NumberFormat = "YYYY-MM-DD HH:MM:SS"
